# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-27 20:13:34
#
# The "Recorded By" column (G) lists the users who recorded/edited an
# attendance entry as a comma-separated string (e.g. "System, user@example.com").
# Upstream reordered these lists (oldest/last recorder now listed first), so for
# every row reverse the order of the comma-separated entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $text = $cell.Value2

    if ($text -ne $null -and $text -is [string] -and $text.Contains(",")) {
        $parts = $text -split ", "
        $n = $parts.Count
        if ($n -gt 1) {
            $reversed = $parts[($n - 1)..0]
            $newText = [string]::Join(", ", $reversed)
            # Use .Equals() (ordinal/case-sensitive) rather than -ne, since the
            # PowerShell comparison operators here are case-insensitive and would
            # otherwise skip rows that only differ by letter casing (e.g. the
            # "System" / "system" duplicate-looking entries).
            if (-not $newText.Equals($text)) {
                $cell.Value = $newText
            }
        }
    }
}
